# Release Aspose.Cells Cloud SDK 23.12 — workbook gains data on Sheet1 plus
# two more (initially blank) worksheets, each with an explicit page setup.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# Sheet1: put the value in A1.
$ws1.Range("A1").Value = 1111

# Add Sheet2 then Sheet3, each appended after the current last sheet so the
# tab order ends up Sheet1, Sheet2, Sheet3.
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets($wb.Worksheets.Count))
$ws2.Name = "Sheet2"
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets($wb.Worksheets.Count))
$ws3.Name = "Sheet3"

# Give every sheet an explicit A4/portrait page setup.
foreach ($sheet in @($ws1, $ws2, $ws3)) {
    $sheet.PageSetup.PaperSize = 9
    $sheet.PageSetup.Orientation = 1
}

# Leave the workbook focused back on Sheet1 with A2 selected, matching the
# view state the workbook was saved with.
$null = $ws1.Activate()
$null = $ws1.Range("A2").Select()
